$d = $word.ActiveDocument

# Locate the paragraph that currently holds "Recall they overlap..." (and the _GoBack bookmark)
$targetIdx = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $pp = $d.Paragraphs.Item($i)
    if ($pp.Range.Text -like "Recall they overlap*") {
        $targetIdx = $i
        break
    }
}

$target = $d.Paragraphs.Item($targetIdx)
$tr = $target.Range

# Insert three new empty paragraphs immediately before it:
#   [empty] [new text paragraph] [empty] [Recall they overlap...]
$tr.InsertParagraphBefore()
$tr.InsertParagraphBefore()
$tr.InsertParagraphBefore()

$newTextParaIdx = $targetIdx + 1
$newPara = $d.Paragraphs.Item($newTextParaIdx)
$newPara.Range.Text = "Randomly sample per each base pair? For distribution "

# Move the _GoBack bookmark from the "Recall..." paragraph onto the end of the
# new "Randomly sample..." paragraph (right after its run, same as it used to sit
# right after the run in the old "Recall..." paragraph).
$bm = $d.Bookmarks.Item("_GoBack")
$bm.Delete()

$npRange = $d.Paragraphs.Item($newTextParaIdx).Range
$npRange.MoveEnd(1, -1) | Out-Null
$endPos = $npRange.End

# Placing a zero-length bookmark exactly at the end of a paragraph's text can make
# it bleed into the start of the following paragraph, so temporarily add a
# placeholder character after the insertion point, anchor the bookmark there, then
# remove the placeholder again.
$placeholderRange = $d.Range($endPos, $endPos)
$placeholderRange.InsertAfter("X")

$bmRange = $d.Range($endPos, $endPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

$delRange = $d.Range($endPos, $endPos + 1)
$delRange.Delete()
